# Traded. Fixed 20 minute trade problem.
# Append the new trade row (row 4) to the HZNP random-trade log.
#
# Columns: A=Date  B=Profitable  C=Principle  D=Start Principle
#          E=BuyPrice  F=SellPrice  G=IsShortSell  H=Price Change %
#          I=Strong trade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (A3 / G3) carries the date-time display style; copy that
# formatting down onto the new row's date / short-sell cells before
# writing their values, so row 4 matches the existing formatting.
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("G3").Copy($ws.Range("G4"))

$ws.Range("A4").Value = 42641.539583333331
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9880.61
$ws.Range("D4").Value = 9891
$ws.Range("E4").Value = 18.670000000000002
$ws.Range("F4").Value = 18.71
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 0.21
$ws.Range("I4").Value = $false
